$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2's values into a new row 3 (same data repeated)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = 0.002
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.1
$ws.Range("M3").Value = 0.01
